$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all Tc (Test Case results) to "Y" (Yes)
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"

# Update the active selection to C2
$ws.Range("C2").Select()
